$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.165.32'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '2.588.34'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = '2.598.46'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -4.55%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("E13").Value = '  +3.50%  '
$ws.Range("D14").Value = '3.045.64'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '58.938.58'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '2.648.69'
$ws.Range("E17").Value = '  +2.45%  '
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '338.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.44%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("E35").Value = '  -1.62%  '
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.818'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.45%  '
$ws.Range("E41").Value = '  -0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '272.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0953'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0516'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("D49").Value = '1.962.63'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0220'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.58%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.71%  '
